$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "number of iterations" column within each of the three
# result blocks (secao_aurea, armijo, interpolacao). Each block originally
# spanned 3 columns (e.g. B:D); a new column is inserted right after the
# first column of each block (e.g. after B -> new col before former E),
# pushing the remaining columns to the right.
$ws.Columns("E:E").Insert()
$ws.Columns("I:I").Insert()

# Fill the newly inserted "iterations" columns (E, I) plus the new trailing
# column (M) with the iteration counts for each method/function row.
$ws.Range("E2").Value = 315
$ws.Range("I2").Value = 2000
$ws.Range("M2").Value = 1475

$ws.Range("E3").Value = 10000
$ws.Range("I3").Value = 10000
$ws.Range("M3").Value = 10000

$ws.Range("E4").Value = 20
$ws.Range("I4").Value = 31
$ws.Range("M4").Value = 9

$ws.Range("E5").Value = 9
$ws.Range("I5").Value = 10
$ws.Range("M5").Value = 6

$ws.Range("E6").Value = 13
$ws.Range("I6").Value = 16
$ws.Range("M6").Value = 85
